$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the first commitment date cell (1/20/2023 => serial 44946) and its date number format
$ws.Range("J2").Value = 44946
$ws.Range("J2").NumberFormat = "mm-dd-yy"

# Propagate the same format to the remaining rows via copy/paste-format so they
# all share a single cell style (rather than each getting its own new style)
$ws.Range("J2").Copy()
$ws.Range("J3:J7").PasteSpecial(-4122)
$ws.Range("J3:J7").Value = 44946
$excel.CutCopyMode = $false

# Column J should best-fit the (date) contents
$ws.Columns.Item(10).AutoFit()

# Add header for new "Commitment Date" column
$ws.Range("J1").Value = "Commitment Date"

# Update selection to mirror the edited workbook state
$ws.Range("J3:J7").Select()
